# Re-execute the RAD "Extension Payments" test row: refresh its execution
# timestamp and clear the stale Execute flags on the other rows so only
# this row is picked up for the next RAD run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Quarterly Estimated Tax) and Row 4 (New Tax Return Amount Due):
# remove their "Execute" flag cell entirely.
$ws.Range("C2").Clear()
$ws.Range("C4").Clear()

# Row 3 (Extension Payments): stamp the new execution date/time.
$ws.Range("B3").Value = "Wed Mar 20 23:05:45 EDT 2024"
